$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (61) down into the three
# new rows (62:64) so the new rows pick up the same cell style (s="1").
$ws.Range("A61:R61").Copy()
$ws.Range("A62:R64").PasteSpecial(-4122)

# Row 62 - new survey response (4/17/2021 14:14:25)
$ws.Range("A62").Value = "4/17/2021 14:14:25"
$ws.Range("C62").Value = "Beginner"
$ws.Range("D62").Value = "Intermediate"
$ws.Range("E62").Value = 2
$ws.Range("F62").Value = "AI"
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = "AI"
$ws.Range("I62").Value = 1
$ws.Range("J62").Value = "AI"
$ws.Range("K62").Value = 2
$ws.Range("L62").Value = "AI"
$ws.Range("M62").Value = 1
$ws.Range("N62").Value = "AI"
$ws.Range("O62").Value = "3 - Neutral"
$ws.Range("P62").Value = 4
$ws.Range("Q62").Value = 2
$ws.Range("R62").Value = 2

# Row 63 - new survey response (4/17/2021 15:02:26)
$ws.Range("A63").Value = "4/17/2021 15:02:26"
$ws.Range("C63").Value = "Intermediate"
$ws.Range("D63").Value = "Intermediate"
$ws.Range("E63").Value = 2
$ws.Range("F63").Value = "AI"
$ws.Range("G63").Value = 3
$ws.Range("H63").Value = "Human"
$ws.Range("I63").Value = 1
$ws.Range("J63").Value = "AI"
$ws.Range("K63").Value = 3
$ws.Range("L63").Value = "Human"
$ws.Range("M63").Value = 1
$ws.Range("N63").Value = "AI"
$ws.Range("O63").Value = 2
$ws.Range("P63").Value = 4
$ws.Range("Q63").Value = "1 - Sad/Dark"
$ws.Range("R63").Value = 2

# Row 64 - new survey response (4/17/2021 15:50:31)
$ws.Range("A64").Value = "4/17/2021 15:50:31"
$ws.Range("C64").Value = "Beginner"
$ws.Range("D64").Value = "Beginner"
$ws.Range("E64").Value = 5
$ws.Range("F64").Value = "Human"
$ws.Range("G64").Value = 5
$ws.Range("H64").Value = "Human"
$ws.Range("I64").Value = 5
$ws.Range("J64").Value = "Human"
$ws.Range("K64").Value = 5
$ws.Range("L64").Value = "Human"
$ws.Range("M64").Value = 5
$ws.Range("N64").Value = "Human"
$ws.Range("O64").Value = "5 - Happy/Bright"
$ws.Range("P64").Value = "5 - Bright/Happy"
$ws.Range("Q64").Value = "5 - Happy/Bright"
$ws.Range("R64").Value = "5 - Bright/Happy"

# Freeze the header row and restore the normal (non-split) pane state, then
# move the active selection to D9 as in the source workbook.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D9").Select()
